# Insert two new price records for "Vega Modelo de Temuco - Mango" dataset.
# This shifts the existing rows 342:460 down to 344:462 and populates the
# two newly-inserted rows (342 and 343) with the new data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 342, pushing old row 342 onward down by two.
$ws.Rows("342:343").Insert()

# --- Row 342 (new record) ---
$ws.Cells.Item(342, 1).Value2  = 10
$ws.Cells.Item(342, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(342, 3).Value2  = "La Araucanía"
$ws.Cells.Item(342, 4).Value2  = 44900
$ws.Cells.Item(342, 5).Value2  = 9
$ws.Cells.Item(342, 6).Value2  = "Fruta"
$ws.Cells.Item(342, 7).Value2  = 100108
$ws.Cells.Item(342, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(342, 9).Value2  = 100108002
$ws.Cells.Item(342, 10).Value2 = "Mango"
$ws.Cells.Item(342, 11).Value2 = "Sin especificar"
$ws.Cells.Item(342, 12).Value2 = "Primera"
$ws.Cells.Item(342, 13).Value2 = 2200
$ws.Cells.Item(342, 14).Value2 = 8000
$ws.Cells.Item(342, 15).Value2 = 9000
$ws.Cells.Item(342, 16).Value2 = 8545
$ws.Cells.Item(342, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(342, 18).Value2 = "Perú"
$ws.Cells.Item(342, 19).Value2 = 2136
$ws.Cells.Item(342, 20).Value2 = 4

# --- Row 343 (new record) ---
$ws.Cells.Item(343, 1).Value2  = 10
$ws.Cells.Item(343, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(343, 3).Value2  = "La Araucanía"
$ws.Cells.Item(343, 4).Value2  = 44900
$ws.Cells.Item(343, 5).Value2  = 9
$ws.Cells.Item(343, 6).Value2  = "Fruta"
$ws.Cells.Item(343, 7).Value2  = 100108
$ws.Cells.Item(343, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(343, 9).Value2  = 100108002
$ws.Cells.Item(343, 10).Value2 = "Mango"
$ws.Cells.Item(343, 11).Value2 = "Sin especificar"
$ws.Cells.Item(343, 12).Value2 = "Segunda"
$ws.Cells.Item(343, 13).Value2 = 500
$ws.Cells.Item(343, 14).Value2 = 7000
$ws.Cells.Item(343, 15).Value2 = 7000
$ws.Cells.Item(343, 16).Value2 = 7000
$ws.Cells.Item(343, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(343, 18).Value2 = "Perú"
$ws.Cells.Item(343, 19).Value2 = 1750
$ws.Cells.Item(343, 20).Value2 = 4
